# edit.ps1 -- apply the "Culinary Crossroads" -> "Marvelous World of Biology"
# content swap described by the diff.
#
# NOTE: we deliberately avoid Find.Execute's built-in Replace (it runs the
# text through smart-quote AutoCorrect, turning straight apostrophes into
# curly ones). Instead we Find (search only) to locate/select the target
# range, then assign Range.Text directly, which substitutes the text
# byte-for-byte and keeps straight apostrophes intact.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "text not found: $old"
    }
    $rng.Text = $new
    return $rng
}

# --- Title / byline / email -------------------------------------------------
Replace-Text "Culinary Crossroads: Exploring Global Connections in Cuisine" "The Marvelous World of Biology: Exploring the Secrets of Life" | Out-Null
Replace-Text "Amelia Gomez, PhD" "Alexandria Hayes" | Out-Null
Replace-Text "agomez@culinaryinstitute" "alexandriahayes@emailworld" | Out-Null
Replace-Text "edu" "net" | Out-Null

# --- Body paragraph ----------------------------------------------------------
Replace-Text "As humans, our connection to food transcends mere sustenance; it weaves a rich tapestry of culture, history, and identity" "Biology, the study of life, embarks on an extraordinary journey into the intricate workings of living organisms" | Out-Null

Replace-Text " Throughout history, cuisine has served as a conduit for global exchange, bridging diverse regions and fostering intercultural dialogue" " As we delve into the diverse tapestry of nature's marvels, we unravel the mysteries that govern the behavior and interactions of plants, animals, and microorganisms" | Out-Null

Replace-Text " From the spice trade that reshaped culinary landscapes to the diaspora of people carrying beloved recipes across borders, food has consistently played a pivotal role in shaping our global community" " From the smallest cell to the vast ecosystems, biology captivates our imagination and challenges us to comprehend the enigmatic dance of life" | Out-Null

Replace-Text "This essay delves into the myriad ways in which cuisine has facilitated global interconnectedness, showcasing the fascinating interplay between food, culture, and history" "Biology serves as a bridge between the grandeur of the universe and the intricacies of our own bodies" | Out-Null

# This sentence gains two extra trailing sentences in the new version.
$rng = Replace-Text " Whether it be through the dissemination of ingredients, the fusion of culinary techniques, or the influence of migration patterns, food has acted as a potent force in fostering understanding and appreciation among people from vastly different backgrounds" " It unravels the secrets of DNA, the blueprint of heredity, and uncovers the mysteries of genetic variation, the driving force of evolution and adaptation"
$ip = $rng.Duplicate
$ip.Collapse(0)
$ip.InsertAfter(". We witness the enchanting spectacle of cells dividing, organizing, and communicating, their ceaseless symphony orchestrating the wonders of life")

Replace-Text "The exploration begins by examining the historical significance of trade routes, particularly the Silk Road, in facilitating the exchange of spices, herbs, and other culinary treasures" "With each passing day, biological discoveries illuminate the path toward medical advancements and environmental solutions" | Out-Null

$rng2 = Replace-Text " This vibrant network of interconnectedness not only introduced novel flavors and ingredients to various regions but also spurred innovation and experimentation within the culinary realm" " We gain invaluable insights into diseases, their causes, and potential treatments, offering hope to those touched by illness"
$ip2 = $rng2.Duplicate
$ip2.Collapse(0)
$ip2.InsertAfter(". As we explore the interconnections within ecosystems, we unravel the intricate web of life's dependencies and strive for sustainable practices that harmonize human activities with the natural world")

# --- Summary paragraph --------------------------------------------------------
Replace-Text "This essay unveils the profound role that cuisine has played in fostering global interconnectedness, demonstrating how food has served as a catalyst for cultural exchange, innovation, and understanding among diverse communities" "Biology, the captivating study of life, unveils the complexities and marvels of living organisms" | Out-Null

$rng3 = Replace-Text " The exploration of historical trade routes, the impact of cultural diffusion, and the influence of migration patterns reveals the remarkable ways in which cuisine has transcended geographic boundaries, enriching our collective culinary heritage and fostering a sense of global kinship" " From the enigmatic dance of cells to the majesty of ecosystems, biology inspires awe and wonder"
$ip3 = $rng3.Duplicate
$ip3.Collapse(0)
$ip3.InsertAfter(". This science uncovers the mysteries of heredity, evolution, and adaptation, illuminating the tapestry of life's astonishing diversity. It plays a pivotal role in medical advancements and environmental solutions, offering hope for a healthier and harmonious world. Biology, a symphony of discovery, invites us to explore the captivating secrets of life, revealing the profound interconnectedness of all living things")

# --- Trailing empty paragraph --------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter() | Out-Null

Write-Output "edit complete"
